$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028662385041881
$ws.Range("D2").Value = 1.034478538279474
$ws.Range("E2").Value = 1.05196914091942
$ws.Range("F2").Value = 1.057278636590748
$ws.Range("I2").Value = 1.037740101213493
$ws.Range("J2").Value = 1.033813459315038
$ws.Range("K2").Value = 1.037277896269257
$ws.Range("L2").Value = 1.054719171552567
$ws.Range("M2").Value = 1.060014043168019
$ws.Range("N2").Value = 1.015298137878576
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029425889528576
$ws.Range("D3").Value = 1.03505984993659
$ws.Range("E3").Value = 1.052997086803344
$ws.Range("F3").Value = 1.058316894026917
$ws.Range("I3").Value = 1.037935144593535
$ws.Range("J3").Value = 1.03421862074644
$ws.Range("K3").Value = 1.037669011358081
$ws.Range("L3").Value = 1.055559316823866
$ws.Range("M3").Value = 1.060865535046
$ws.Range("N3").Value = 1.01543216461023
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029920251642668
$ws.Range("D4").Value = 1.035436171099503
$ws.Range("E4").Value = 1.053663441769058
$ws.Range("F4").Value = 1.058989726098622
$ws.Range("I4").Value = 1.038060147304765
$ws.Range("J4").Value = 1.034480427248987
$ws.Range("K4").Value = 1.03792156249964
$ws.Range("L4").Value = 1.056103543076214
$ws.Range("M4").Value = 1.061416909404649
$ws.Range("N4").Value = 1.015518759869426
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030128157355732
$ws.Range("D5").Value = 1.035594416410611
$ws.Range("E5").Value = 1.053943864853247
$ws.Range("F5").Value = 1.059272825346036
$ws.Range("I5").Value = 1.038112409604061
$ws.Range("J5").Value = 1.03459040342484
$ws.Range("K5").Value = 1.038027607796236
$ws.Range("L5").Value = 1.056332477528944
$ws.Range("M5").Value = 1.06164880238642
$ws.Range("N5").Value = 1.015555133251548
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030163070042817
$ws.Range("D6").Value = 1.035620988785213
$ws.Range("E6").Value = 1.053990965948302
$ws.Range("F6").Value = 1.059320373032938
$ws.Range("I6").Value = 1.038121167714694
$ws.Range("J6").Value = 1.034608863759122
$ws.Range("K6").Value = 1.03804540577501
$ws.Range("L6").Value = 1.056370924914388
$ws.Range("M6").Value = 1.061687743789788
$ws.Range("N6").Value = 1.015561238657943
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029923029392832
$ws.Range("D7").Value = 1.035438285426452
$ws.Range("E7").Value = 1.053667187667869
$ws.Range("F7").Value = 1.058993507940512
$ws.Range("I7").Value = 1.038060846772042
$ws.Range("J7").Value = 1.034481897100334
$ws.Range("K7").Value = 1.037922979983368
$ws.Range("L7").Value = 1.056106601553853
$ws.Range("M7").Value = 1.061420007597291
$ws.Range("N7").Value = 1.015519246015668
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028920346978269
$ws.Range("D8").Value = 1.034674958312834
$ws.Range("E8").Value = 1.052316289838975
$ws.Range("F8").Value = 1.057629310701756
$ws.Range("I8").Value = 1.037806265760656
$ws.Range("J8").Value = 1.03395045941388
$ws.Range("K8").Value = 1.037410183755094
$ws.Range("L8").Value = 1.055002978574733
$ws.Range("M8").Value = 1.060301724910806
$ws.Range("N8").Value = 1.015343459320273
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027156042138697
$ws.Range("D9").Value = 1.033331282379326
$ws.Range("E9").Value = 1.049945113818833
$ws.Range("F9").Value = 1.055233208755727
$ws.Range("I9").Value = 1.037348477798091
$ws.Range("J9").Value = 1.033011288733924
$ws.Range("K9").Value = 1.036502590504434
$ws.Range("L9").Value = 1.053062861971313
$ws.Range("M9").Value = 1.058334295926198
$ws.Range("N9").Value = 1.015032729603148
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02598164833908
$ws.Range("D10").Value = 1.032436542077265
$ws.Range("E10").Value = 1.048370631658418
$ws.Range("F10").Value = 1.053641112458004
$ws.Range("I10").Value = 1.037037155383951
$ws.Range("J10").Value = 1.032383417494929
$ws.Range("K10").Value = 1.03589492096805
$ws.Range("L10").Value = 1.051772605966631
$ws.Range("M10").Value = 1.057024843260303
$ws.Range("N10").Value = 1.014824946225199
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025473571982983
$ws.Range("D11").Value = 1.032049378197583
$ws.Range("E11").Value = 1.047690371268232
$ws.Range("F11").Value = 1.052952989639804
$ws.Range("I11").Value = 1.036900904517431
$ws.Range("J11").Value = 1.032111138960999
$ws.Range("K11").Value = 1.035631189419329
$ws.Range("L11").Value = 1.051214670908319
$ws.Range("M11").Value = 1.05645836281708
$ws.Range("N11").Value = 1.014734829006183
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025284918088044
$ws.Range("D12").Value = 1.031905609737097
$ws.Range("E12").Value = 1.047437918760294
$ws.Range("F12").Value = 1.052697581399214
$ws.Range("I12").Value = 1.036850078322574
$ws.Range("J12").Value = 1.032009942644636
$ws.Range("K12").Value = 1.035533137794587
$ws.Range("L12").Value = 1.05100754321742
$ws.Range("M12").Value = 1.056248026080612
$ws.Range("N12").Value = 1.014701333956373
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025325381873247
$ws.Range("D13").Value = 1.031936446657385
$ws.Range("E13").Value = 1.047492060414729
$ws.Range("F13").Value = 1.052752358663275
$ws.Range("I13").Value = 1.036860990501978
$ws.Range("J13").Value = 1.032031652293468
$ws.Range("K13").Value = 1.035554174254429
$ws.Range("L13").Value = 1.0510519676264
$ws.Range("M13").Value = 1.056293140418673
$ws.Range("N13").Value = 1.014708519725194
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02545797639095
$ws.Range("D14").Value = 1.032037493396591
$ws.Range("E14").Value = 1.047669498838961
$ws.Range("F14").Value = 1.052931873607879
$ws.Range("I14").Value = 1.036896707621829
$ws.Range("J14").Value = 1.032102775258592
$ws.Range("K14").Value = 1.035623086270435
$ws.Range("L14").Value = 1.051197547333221
$ws.Range("M14").Value = 1.056440974692762
$ws.Range("N14").Value = 1.014732060731363
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02553968134568
$ws.Range("D15").Value = 1.032099757169741
$ws.Range("E15").Value = 1.047778854550873
$ws.Range("F15").Value = 1.053042504060633
$ws.Range("I15").Value = 1.036918685443894
$ws.Range("J15").Value = 1.032146588553452
$ws.Range("K15").Value = 1.035665533348507
$ws.Range("L15").Value = 1.051287258959346
$ws.Range("M15").Value = 1.056532070809742
$ws.Range("N15").Value = 1.014746562284301
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026015377182192
$ws.Range("D16").Value = 1.03246224258782
$ws.Range("E16").Value = 1.048415810006246
$ws.Range("F16").Value = 1.053686807657389
$ws.Range("I16").Value = 1.037046167496736
$ws.Range("J16").Value = 1.032401479259113
$ws.Range("K16").Value = 1.035912411279942
$ws.Range("L16").Value = 1.051809650248498
$ws.Range("M16").Value = 1.057062449769564
$ws.Range("N16").Value = 1.014830923972049
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026313888584629
$ws.Range("D17").Value = 1.032689692145073
$ws.Range("E17").Value = 1.04881575785016
$ws.Range("F17").Value = 1.054091301837692
$ws.Range("I17").Value = 1.037125746937718
$ws.Range("J17").Value = 1.032561257486071
$ws.Range("K17").Value = 1.03606710968736
$ws.Range("L17").Value = 1.052137535089807
$ws.Range("M17").Value = 1.057395282736309
$ws.Range("N17").Value = 1.014883803104599
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026488047959173
$ws.Range("D18").Value = 1.032822385038059
$ws.Range("E18").Value = 1.049049185395585
$ws.Range("F18").Value = 1.05432735846657
$ws.Range("I18").Value = 1.037172024692976
$ws.Range("J18").Value = 1.032654414163714
$ws.Range("K18").Value = 1.036157284042852
$ws.Range("L18").Value = 1.052328857488927
$ws.Range("M18").Value = 1.057589468794621
$ws.Range("N18").Value = 1.014914632529915
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026547439003448
$ws.Range("D19").Value = 1.032867634144609
$ws.Range("E19").Value = 1.049128802674982
$ws.Range("F19").Value = 1.0544078683395
$ws.Range("I19").Value = 1.037187780513451
$ws.Range("J19").Value = 1.032686171491756
$ws.Range("K19").Value = 1.036188021191639
$ws.Range("L19").Value = 1.052394105789663
$ws.Range("M19").Value = 1.057655689709563
$ws.Range("N19").Value = 1.014925142167811
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026281856705585
$ws.Range("D20").Value = 1.032665286319181
$ws.Range("E20").Value = 1.048772832253633
$ws.Range("F20").Value = 1.054047890811453
$ws.Range("I20").Value = 1.037117223249697
$ws.Range("J20").Value = 1.032544118847425
$ws.Range("K20").Value = 1.036050518066083
$ws.Range("L20").Value = 1.052102348621217
$ws.Range("M20").Value = 1.05735956770702
$ws.Range("N20").Value = 1.014878131127537
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025418928712497
$ws.Range("D21").Value = 1.032007736492814
$ws.Range("E21").Value = 1.047617241385159
$ws.Range("F21").Value = 1.05287900564369
$ws.Range("I21").Value = 1.036886195792494
$ws.Range("J21").Value = 1.032081832964239
$ws.Range("K21").Value = 1.035602795874069
$ws.Range("L21").Value = 1.051154674577383
$ws.Range("M21").Value = 1.056397438985821
$ws.Range("N21").Value = 1.014725129082313
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024876767503388
$ws.Range("D22").Value = 1.031594549429759
$ws.Range("E22").Value = 1.046891986446608
$ws.Range("F22").Value = 1.052145187994199
$ws.Range("I22").Value = 1.036739687088626
$ws.Range("J22").Value = 1.031790829184688
$ws.Range("K22").Value = 1.035320775107291
$ws.Range("L22").Value = 1.050559495655386
$ws.Range("M22").Value = 1.055792969916354
$ws.Range("N22").Value = 1.01462880640611
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025164139362338
$ws.Range("D23").Value = 1.031813564338227
$ws.Range("E23").Value = 1.047276333212574
$ws.Range("F23").Value = 1.052534093328197
$ws.Range("I23").Value = 1.036817472613411
$ws.Range("J23").Value = 1.031945128261135
$ws.Range("K23").Value = 1.035470328561912
$ws.Range("L23").Value = 1.050874948269695
$ws.Range("M23").Value = 1.056113366515734
$ws.Range("N23").Value = 1.014679880527145
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026296330401105
$ws.Range("D24").Value = 1.03267631418305
$ws.Range("E24").Value = 1.048792228035857
$ws.Range("F24").Value = 1.054067506008671
$ws.Range("I24").Value = 1.037121075168791
$ws.Range("J24").Value = 1.032551863181933
$ws.Range("K24").Value = 1.036058015286785
$ws.Range("L24").Value = 1.052118247647162
$ws.Range("M24").Value = 1.057375705636323
$ws.Range("N24").Value = 1.014880694093367
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027611845273509
$ws.Range("D25").Value = 1.033678478148758
$ws.Range("E25").Value = 1.050557013580739
$ws.Range("F25").Value = 1.055851728562964
$ws.Range("I25").Value = 1.037467910726562
$ws.Range("J25").Value = 1.03325440133421
$ws.Range("K25").Value = 1.036737689733027
$ws.Range("L25").Value = 1.053563876696371
$ws.Range("M25").Value = 1.05884254626203
$ws.Range("N25").Value = 1.015113173424346

Write-Host "updated 264 cells"